$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2528.5
$ws.Cells.Item(17, 10).Value = 2528.5
$ws.Cells.Item(17, 12).Value = 7585.5
$ws.Cells.Item(17, 14).Value = -7921.5
$ws.Cells.Item(31, 8).Value = 232.16667
$ws.Cells.Item(31, 9).Value = 232.16667
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 696.50001
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -466.50001
$ws.Cells.Item(31, 14).ClearContents()
$ws.Cells.Item(98, 8).Value = 1459.4
$ws.Cells.Item(98, 9).Value = 627.7143
$ws.Cells.Item(98, 10).Value = 3400
$ws.Cells.Item(98, 11).Value = 627.7143
$ws.Cells.Item(98, 12).Value = 3400
$ws.Cells.Item(98, 13).Value = 870.2857
$ws.Cells.Item(98, 14).Value = -6396
$ws.Cells.Item(122, 8).Value = 1459.4
$ws.Cells.Item(122, 9).Value = 627.7143
$ws.Cells.Item(122, 10).Value = 3400
$ws.Cells.Item(122, 11).Value = 1883.1429
$ws.Cells.Item(122, 12).Value = 10200
$ws.Cells.Item(122, 13).Value = 566.8571000000002
$ws.Cells.Item(122, 14).Value = -15100
$ws.Cells.Item(123, 8).Value = 70777
$ws.Cells.Item(123, 10).Value = 70777
$ws.Cells.Item(123, 12).Value = 70777
$ws.Cells.Item(123, 14).Value = -80577
$ws.Cells.Item(129, 8).Value = 4056.9
$ws.Cells.Item(129, 9).Value = 4109.857
$ws.Cells.Item(129, 10).Value = 3933.3333
$ws.Cells.Item(129, 11).Value = 12329.571
$ws.Cells.Item(129, 12).Value = 11799.9999
$ws.Cells.Item(129, 13).Value = -7329.571
$ws.Cells.Item(129, 14).Value = -21799.9999
$ws.Cells.Item(132, 8).Value = 233485.05
$ws.Cells.Item(132, 9).Value = 948.4872
$ws.Cells.Item(132, 11).Value = 2845.4616
$ws.Cells.Item(132, 13).Value = -315.4616000000001
$ws.Cells.Item(137, 8).Value = 3305.6365
$ws.Cells.Item(137, 10).Value = 3530.5
$ws.Cells.Item(137, 12).Value = 10591.5
$ws.Cells.Item(137, 14).Value = -15691.5
$ws.Cells.Item(138, 8).Value = 2109.2551
$ws.Cells.Item(138, 9).Value = 860.2222
$ws.Cells.Item(138, 10).Value = 2584.2395
$ws.Cells.Item(138, 11).Value = 2580.6666
$ws.Cells.Item(138, 12).Value = 7752.718500000001
$ws.Cells.Item(138, 13).Value = 2559.3334
$ws.Cells.Item(138, 14).Value = -18032.7185
$ws.Cells.Item(141, 8).Value = 3578.6843
$ws.Cells.Item(141, 9).Value = 2187.375
$ws.Cells.Item(141, 11).Value = 6562.125
$ws.Cells.Item(141, 13).Value = -1382.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2335.4648
$ws.Cells.Item(32, 9).Value = 1263.7462
$ws.Cells.Item(32, 11).Value = 1263.7462
$ws.Cells.Item(32, 13).Value = -976.7462
$ws.Cells.Item(39, 8).Value = 4800
$ws.Cells.Item(39, 9).Value = 4800
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 4800
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = -4280
$ws.Cells.Item(39, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 3456.0435
$ws.Cells.Item(61, 9).Value = 3239.95
$ws.Cells.Item(61, 10).Value = 4896.6665
$ws.Cells.Item(61, 11).Value = 3239.95
$ws.Cells.Item(61, 12).Value = 4896.6665
$ws.Cells.Item(61, 13).Value = -3027.95
$ws.Cells.Item(61, 14).Value = -5320.6665
$ws.Cells.Item(74, 8).Value = 1131.7858
$ws.Cells.Item(74, 9).Value = 1164.091
$ws.Cells.Item(74, 11).Value = 1164.091
$ws.Cells.Item(74, 13).Value = -290.0909999999999
$ws.Cells.Item(77, 8).Value = 1131.7858
$ws.Cells.Item(77, 9).Value = 1164.091
$ws.Cells.Item(77, 11).Value = 5820.455
$ws.Cells.Item(77, 13).Value = -1452.455
$ws.Cells.Item(122, 8).Value = 3176.3044
$ws.Cells.Item(122, 9).Value = 3239.0588
$ws.Cells.Item(122, 10).Value = 2998.5
$ws.Cells.Item(122, 11).Value = 9717.1764
$ws.Cells.Item(122, 12).Value = 8995.5
$ws.Cells.Item(122, 13).Value = -7267.1764
$ws.Cells.Item(122, 14).Value = -13895.5
$ws.Cells.Item(132, 8).Value = 1616.575
$ws.Cells.Item(132, 9).Value = 1555.4615
$ws.Cells.Item(132, 11).Value = 4666.3845
$ws.Cells.Item(132, 13).Value = -2136.3845
$ws.Cells.Item(136, 8).Value = 3456.0435
$ws.Cells.Item(136, 9).Value = 3239.95
$ws.Cells.Item(136, 10).Value = 4896.6665
$ws.Cells.Item(136, 11).Value = 9719.849999999999
$ws.Cells.Item(136, 12).Value = 14689.9995
$ws.Cells.Item(136, 13).Value = -7169.849999999999
$ws.Cells.Item(136, 14).Value = -19789.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5899.2
$ws.Cells.Item(20, 9).Value = 3222.875
$ws.Cells.Item(20, 11).Value = 3222.875
$ws.Cells.Item(20, 13).Value = -2975.875
$ws.Cells.Item(94, 8).Value = 1231.3478
$ws.Cells.Item(94, 9).Value = 654.5294
$ws.Cells.Item(94, 10).Value = 2865.6667
$ws.Cells.Item(94, 11).Value = 654.5294
$ws.Cells.Item(94, 12).Value = 2865.6667
$ws.Cells.Item(94, 13).Value = -203.5294
$ws.Cells.Item(94, 14).Value = -3767.6667
$ws.Cells.Item(105, 8).Value = 2801.7
$ws.Cells.Item(105, 10).Value = 3582.3333
$ws.Cells.Item(105, 12).Value = 3582.3333
$ws.Cells.Item(105, 14).Value = -7076.3333
$ws.Cells.Item(134, 8).Value = 1015.44
$ws.Cells.Item(134, 9).Value = 773.3043
$ws.Cells.Item(134, 11).Value = 2319.9129
$ws.Cells.Item(134, 13).Value = 215.0870999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2539.8
$ws.Cells.Item(34, 8).Value = 2539.8
$ws.Cells.Item(58, 8).Value = 2194.6667
$ws.Cells.Item(58, 9).Value = 2039.7727
$ws.Cells.Item(58, 10).Value = 3898.5
$ws.Cells.Item(58, 11).Value = 2039.7727
$ws.Cells.Item(58, 12).Value = 3898.5
$ws.Cells.Item(58, 13).Value = -1836.7727
$ws.Cells.Item(58, 14).Value = -4304.5
$ws.Cells.Item(59, 8).Value = 45412.332
$ws.Cells.Item(59, 9).Value = 38999
$ws.Cells.Item(59, 10).Value = 46695
$ws.Cells.Item(59, 11).Value = 38999
$ws.Cells.Item(59, 12).Value = 46695
$ws.Cells.Item(59, 13).Value = -37854
$ws.Cells.Item(59, 14).Value = -48985
$ws.Cells.Item(136, 8).Value = 2194.6667
$ws.Cells.Item(136, 9).Value = 2039.7727
$ws.Cells.Item(136, 10).Value = 3898.5
$ws.Cells.Item(136, 11).Value = 6119.3181
$ws.Cells.Item(136, 12).Value = 11695.5
$ws.Cells.Item(136, 13).Value = -3569.3181
$ws.Cells.Item(136, 14).Value = -16795.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 13134307
$ws.Cells.Item(4, 9).Value = 2528445.2
$ws.Cells.Item(4, 11).Value = 7585335.600000001
$ws.Cells.Item(4, 13).Value = -7585223.600000001
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(46, 8).Value = 735067
$ws.Cells.Item(46, 10).Value = 250235
$ws.Cells.Item(46, 12).Value = 750705
$ws.Cells.Item(46, 14).Value = -750887
$ws.Cells.Item(92, 8).Value = 277.81818
$ws.Cells.Item(92, 9).Value = 252.5
$ws.Cells.Item(92, 11).Value = 757.5
$ws.Cells.Item(92, 13).Value = 490.5
$ws.Cells.Item(110, 8).Value = 10538.8
$ws.Cells.Item(110, 9).Value = 10538.8
$ws.Cells.Item(110, 11).Value = 31616.4
$ws.Cells.Item(110, 13).Value = -27526.4
$ws.Cells.Item(129, 8).Value = 2158.5356
$ws.Cells.Item(129, 9).Value = 717.2727
$ws.Cells.Item(129, 10).Value = 3091.1177
$ws.Cells.Item(129, 11).Value = 2151.8181
$ws.Cells.Item(129, 12).Value = 9273.3531
$ws.Cells.Item(129, 13).Value = 2848.1819
$ws.Cells.Item(129, 14).Value = -19273.3531
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 13).ClearContents()
$ws.Cells.Item(139, 8).Value = 9446.154
$ws.Cells.Item(139, 9).Value = 6313.4
$ws.Cells.Item(139, 11).Value = 18940.2
$ws.Cells.Item(139, 13).Value = -13800.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 63874.75
$ws.Cells.Item(39, 10).Value = 63874.75
$ws.Cells.Item(39, 12).Value = 63874.75
$ws.Cells.Item(39, 14).Value = -64938.75
$ws.Cells.Item(80, 8).Value = 13113.238
$ws.Cells.Item(80, 9).Value = 5324.4443
$ws.Cells.Item(80, 10).Value = 18954.834
$ws.Cells.Item(80, 11).Value = 5324.4443
$ws.Cells.Item(80, 12).Value = 18954.834
$ws.Cells.Item(80, 13).Value = -4326.4443
$ws.Cells.Item(80, 14).Value = -20950.834
$ws.Cells.Item(83, 8).Value = 13113.238
$ws.Cells.Item(83, 9).Value = 5324.4443
$ws.Cells.Item(83, 10).Value = 18954.834
$ws.Cells.Item(83, 11).Value = 26622.2215
$ws.Cells.Item(83, 12).Value = 94774.17
$ws.Cells.Item(83, 13).Value = -21630.2215
$ws.Cells.Item(83, 14).Value = -104758.17
$ws.Cells.Item(102, 8).Value = 4557.1113
$ws.Cells.Item(102, 9).Value = 2956.25
$ws.Cells.Item(102, 11).Value = 2956.25
$ws.Cells.Item(102, 13).Value = -1334.25
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 2469.5625
$ws.Cells.Item(122, 9).Value = 2347.2307
$ws.Cells.Item(122, 11).Value = 7041.6921
$ws.Cells.Item(122, 13).Value = -4591.6921
$ws.Cells.Item(132, 8).Value = 1396.4193
$ws.Cells.Item(132, 9).Value = 885.3333
$ws.Cells.Item(132, 10).Value = 3148.7144
$ws.Cells.Item(132, 11).Value = 2655.9999
$ws.Cells.Item(132, 12).Value = 9446.143199999999
$ws.Cells.Item(132, 13).Value = -125.9998999999998
$ws.Cells.Item(132, 14).Value = -14506.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 576.5
$ws.Cells.Item(82, 9).Value = 593.6667
$ws.Cells.Item(82, 11).Value = 593.6667
$ws.Cells.Item(82, 13).Value = -232.6667
$ws.Cells.Item(85, 8).Value = 576.5
$ws.Cells.Item(85, 9).Value = 593.6667
$ws.Cells.Item(85, 11).Value = 593.6667
$ws.Cells.Item(85, 13).Value = 654.3333
$ws.Cells.Item(136, 8).Value = 2618.3809
$ws.Cells.Item(136, 9).Value = 2817.125
$ws.Cells.Item(136, 10).Value = 1982.4
$ws.Cells.Item(136, 11).Value = 8451.375
$ws.Cells.Item(136, 12).Value = 5947.200000000001
$ws.Cells.Item(136, 13).Value = -5901.375
$ws.Cells.Item(136, 14).Value = -11047.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(133, 8).Value = 67500
$ws.Cells.Item(133, 10).Value = 67500
$ws.Cells.Item(133, 12).Value = 67500
$ws.Cells.Item(133, 14).Value = -77620
$ws.Cells.Item(136, 8).Value = 1508.5333
$ws.Cells.Item(136, 10).Value = 1642
$ws.Cells.Item(136, 12).Value = 4926
$ws.Cells.Item(136, 14).Value = -10026

Write-Output "done"